# Add two new mission rows to the Rottenburg sheet (row 7: "Calendartillery"
# marked as "int" difficulty, row 8: "Germanium Gearbox" marked as "exp"
# difficulty), widen column B slightly so the longer names fit, and leave
# the selection where the author left it after typing the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: Calendartillery, flagged in the "int" (C) column.
$ws.Range("B7").Value = "Calendartillery"
$ws.Range("C7").Value = 1

# New row 8: Germanium Gearbox, flagged in the "exp" (E) column.
$ws.Range("B8").Value = "Germanium Gearbox"
$ws.Range("E8").Value = 1

# Column B needs to be a bit wider to fit the new, longer mission names.
$ws.Columns("B").ColumnWidth = 16.86

# Selection ends up on C18 after the edits (matches the saved workbook view).
$ws.Range("C18").Select() | Out-Null
